$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-recognised as a number by
# Excel (single decimal point, no thousands separators). These need the
# cell's number format forced to Text ("@") before the value is written so
# the stored cell stays a literal string, matching the source data (which
# keeps these as plain text, e.g. prices like "590.92").
$textForcedUpdates = @(
    @{ Cell = "D5"; Value = "590.92" }
    @{ Cell = "D6"; Value = "144.78" }
    @{ Cell = "D10"; Value = "6.71" }
    @{ Cell = "D12"; Value = "0.445" }
    @{ Cell = "D13"; Value = "0.0000227" }
    @{ Cell = "D14"; Value = "33.61" }
    @{ Cell = "D18"; Value = "6.73" }
    @{ Cell = "D20"; Value = "430.61" }
    @{ Cell = "D21"; Value = "13.54" }
    @{ Cell = "D23"; Value = "7.11" }
    @{ Cell = "D24"; Value = "81.41" }
    @{ Cell = "D25"; Value = "10.82" }
    @{ Cell = "D27"; Value = "12.01" }
    @{ Cell = "D29"; Value = "2.29" }
    @{ Cell = "D30"; Value = "0.999" }
    @{ Cell = "D33"; Value = "26.69" }
    @{ Cell = "D37"; Value = "5.62" }
    @{ Cell = "D38"; Value = "3.01" }
    @{ Cell = "D39"; Value = "49.65" }
    @{ Cell = "D40"; Value = "0.125" }
    @{ Cell = "D42"; Value = "8.63" }
    @{ Cell = "D44"; Value = "40.18" }
    @{ Cell = "D46"; Value = "375.17" }
    @{ Cell = "D48"; Value = "129.71" }
    @{ Cell = "D50"; Value = "24.12" }
    @{ Cell = "D51"; Value = "0.107" }
)

# Cells whose new text is already unambiguous as text (multi-dot "."-grouped
# numbers, percentage strings with padding, etc.) - Excel keeps these as
# text automatically, so no NumberFormat change is needed/applied.
$plainUpdates = @(
    @{ Cell = "D2"; Value = "60.913.02" }
    @{ Cell = "E2"; Value = "  -3.40%  " }
    @{ Cell = "D3"; Value = "2.910.47" }
    @{ Cell = "E3"; Value = "  -3.90%  " }
    @{ Cell = "E4"; Value = "  -0.04%  " }
    @{ Cell = "E5"; Value = "  -0.94%  " }
    @{ Cell = "E6"; Value = "  -5.43%  " }
    @{ Cell = "E7"; Value = "  -0.02%  " }
    @{ Cell = "E8"; Value = "  -1.55%  " }
    @{ Cell = "D9"; Value = "2.909.79" }
    @{ Cell = "E9"; Value = "  -3.84%  " }
    @{ Cell = "E10"; Value = "  -4.72%  " }
    @{ Cell = "E11"; Value = "  -4.03%  " }
    @{ Cell = "E12"; Value = "  -4.21%  " }
    @{ Cell = "E13"; Value = "  -2.79%  " }
    @{ Cell = "E14"; Value = "  -6.10%  " }
    @{ Cell = "E15"; Value = "  +0.22%  " }
    @{ Cell = "D16"; Value = "3.393.15" }
    @{ Cell = "E16"; Value = "  -3.93%  " }
    @{ Cell = "D17"; Value = "60.848.61" }
    @{ Cell = "E17"; Value = "  -3.45%  " }
    @{ Cell = "E18"; Value = "  -4.97%  " }
    @{ Cell = "D19"; Value = "2.900.72" }
    @{ Cell = "E19"; Value = "  -4.16%  " }
    @{ Cell = "E20"; Value = "  -4.27%  " }
    @{ Cell = "E21"; Value = "  -4.94%  " }
    @{ Cell = "E22"; Value = "  -2.05%  " }
    @{ Cell = "E23"; Value = "  -5.77%  " }
    @{ Cell = "E24"; Value = "  -1.82%  " }
    @{ Cell = "E25"; Value = "  -5.35%  " }
    @{ Cell = "E26"; Value = "  -3.63%  " }
    @{ Cell = "E27"; Value = "  -3.02%  " }
    @{ Cell = "E28"; Value = "  +0.02%  " }
    @{ Cell = "E29"; Value = "  +0.31%  " }
    @{ Cell = "E30"; Value = "  +0.06%  " }
    @{ Cell = "E31"; Value = "  -2.72%  " }
    @{ Cell = "E32"; Value = "  -5.84%  " }
    @{ Cell = "E33"; Value = "  -3.70%  " }
    @{ Cell = "E34"; Value = "  -3.09%  " }
    @{ Cell = "D35"; Value = "0.0₃0859" }
    @{ Cell = "E35"; Value = "  -2.10%  " }
    @{ Cell = "E36"; Value = "  -3.42%  " }
    @{ Cell = "E38"; Value = "  -4.18%  " }
    @{ Cell = "E39"; Value = "  -1.85%  " }
    @{ Cell = "E40"; Value = "  -4.17%  " }
    @{ Cell = "E41"; Value = "  -4.84%  " }
    @{ Cell = "E42"; Value = "  -4.48%  " }
    @{ Cell = "E43"; Value = "  -5.02%  " }
    @{ Cell = "E44"; Value = "  -10.36%  " }
    @{ Cell = "E45"; Value = "  -3.27%  " }
    @{ Cell = "E46"; Value = "  -4.31%  " }
    @{ Cell = "D47"; Value = "2.703.55" }
    @{ Cell = "E47"; Value = "  -0.26%  " }
    @{ Cell = "E48"; Value = "  -3.21%  " }
    @{ Cell = "E49"; Value = "  +0.03%  " }
    @{ Cell = "E50"; Value = "  -10.48%  " }
    @{ Cell = "E51"; Value = "  -2.32%  " }
)

foreach ($u in $textForcedUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

foreach ($u in $plainUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}
